# MEGATAB_EMPREEND_JUN2025vlight.xlsx - "Add files via upload" edit
#
# Repurposes the (empty, header-only) "VER NO MAPA" column E into a combined
# "COORDENADA(DEC)" lat,long column, populated from the existing separate
# LATITUDE (col N) / LONGITUDE (col O) columns, then drops those two columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combined "lat,long" value for each data row (2-23), built from the existing
# LATITUDE (col N) / LONGITUDE (col O) cell values before they are removed.
$coords = @{
    2  = "-3.891234,-38.455678"
    3  = "-3.789012,-38.512345"
    4  = "-3.715678,-38.567890"
    5  = "-3.812345,-38.543210"
    6  = "-3.812345,-38.543210"
    7  = "-3.732456,-38.489123"
    8  = "-3.812345,-38.543210"
    9  = "-3.732456,-38.489123"
    10 = "-3.812345,-38.543210"
    11 = "-3.824532,-38.579120"
    12 = "-3.793874,-38.481262"
    13 = "-3.816781,-38.551234"
    14 = "-3.727890,-38.639012"
    15 = "-3.718389,-38.482273"
    16 = "-3.830000,-38.550000"
    17 = "-3.714701,-38.581138"
    18 = "-3.837602,-38.460851"
    19 = "-3.812563,-38.537415"
    20 = "-3.732028,-38.462216"
    21 = "-3.892758,-38.455388"
    22 = "-3.873219,-38.635111"
    23 = "-3.727890,-38.639012"
}

# 1. Rename the header and fill in the combined coordinate string per row.
$ws.Range("E1").Value = "COORDENADA(DEC)"
foreach ($r in $coords.Keys) {
    $ws.Range("E$r").Value = $coords[$r]
}

# 2. Delete the now-redundant LATITUDE / LONGITUDE columns (N, O).
$ws.Range("N1:O1").EntireColumn.Delete()

# 3. Fix up the AutoFilter so it spans only the live data range again, and
#    repoint the hidden _FilterDatabase defined name at the same range.
$ws.AutoFilterMode = $false
$ws.Range("A1:M23").AutoFilter()
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Planilha1!_FilterDatabase") {
        $n.RefersTo = "=Planilha1!`$A`$1:`$M`$23"
    }
}

# 4. Update sheet view: scroll so column C is leftmost and select the now
#    trailing (empty) columns L:M as the author's saved selection.
$ws.Range("L1:M1048576").Select()
$ws.Application.ActiveWindow.ScrollColumn = 3
